# Update crypto price (D) / volume-change (E) cells per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.252.58"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.906.47"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5250"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07286"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9043"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08184"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "1.508.46"
$ws.Range("E15").Value = "  -20.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008664"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "27.285.14"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.505"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.344"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.739"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.847"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.850"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09251"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8359"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.00%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.984"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.355"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.736"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5771"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.079"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.096"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.568"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4918"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06057"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.72%  "
